$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the new "2509" period, right below the existing
# "2508" data row (row 16), copying its look & feel.
$ws.Rows("17:17").Insert()
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4104)

# Fill the new row with the same worker, new period "2509"
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "7931662"
$ws.Range("D17").Value = "ALBERTO RAMON VASQUEZ SUAREZ"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Update the summary figures: now 2 periods in arrears, total doubled
$ws.Range("E11").Value = 113880
$ws.Range("F13").Value = 2

$excel.CutCopyMode = 0
